# Update "想去人数" (want-to-go count) figures that changed between
# the gh-pages data generation runs (commit 456a3b4).
#
# Sheet "展览" (Exhibitions):
#   F3: 204 -> 205
#   F4: 813 -> 816
#
# Sheet "全部类型" (All Types, combined list):
#   F4: 204 -> 205
#   F5: 813 -> 816

$wb = $excel.ActiveWorkbook

$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F3").Value = 205
$wsExhibit.Range("F4").Value = 816

$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F4").Value = 205
$wsAll.Range("F5").Value = 816
